$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.020.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.69%  '
$ws.Range("D3").Value = '''1.651.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.48%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''215.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.63%  '
$ws.Range("E6").Value = '  +1.61%  '
$ws.Range("D8").Value = '''0.249'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.81%  '
$ws.Range("E9").Value = '  +1.69%  '
$ws.Range("D10").Value = '''19.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("D12").Value = '''1.885.83'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.56%  '
$ws.Range("D13").Value = '''1.649.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.32%  '
$ws.Range("D14").Value = '''4.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.31%  '
$ws.Range("D15").Value = '''0.517'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.84%  '
$ws.Range("D16").Value = '''65.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.93%  '
$ws.Range("D17").Value = '''239.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.03%  '
$ws.Range("D18").Value = '''27.006.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.72%  '
$ws.Range("D19").Value = '''7.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.62%  '
$ws.Range("E20").Value = '  +1.19%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = '''4.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.61%  '
$ws.Range("D23").Value = '''2.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.65%  '
$ws.Range("E24").Value = '  +3.52%  '
$ws.Range("D25").Value = '''145.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.47%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '''7.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.08%  '
$ws.Range("E28").Value = '  +1.86%  '
$ws.Range("E29").Value = '  +2.81%  '
$ws.Range("D30").Value = '''0.0497'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("E31").Value = '  +1.97%  '
$ws.Range("E32").Value = '  +3.43%  '
$ws.Range("D33").Value = '''1.517.68'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("E34").Value = '  +5.23%  '
$ws.Range("E35").Value = '  +8.69%  '
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("E37").Value = '  +2.18%  '
$ws.Range("E38").Value = '  +3.04%  '
$ws.Range("D39").Value = '''0.884'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.39%  '
$ws.Range("E40").Value = '  +2.91%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  +4.28%  '
$ws.Range("D43").Value = '''66.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.37%  '
$ws.Range("D44").Value = '''1.792.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.36%  '
$ws.Range("D45").Value = '''0.773'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.92%  '
$ws.Range("D46").Value = '''0.914'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.40%  '
$ws.Range("D47").Value = '''89.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.19%  '
$ws.Range("D48").Value = '''0.0₆0106'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.96%  '
$ws.Range("E49").Value = '  +2.85%  '
$ws.Range("E50").Value = '  +1.31%  '
$ws.Range("D51").Value = '''0.0977'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.26%  '
